# Revert back to previous pptx:
#  - move the 6 "graph" pictures on slide 1 down slightly so they line up
#    with their bounding boxes (same X, Top += ~5.98pt / 75943 EMU)
#  - refresh the cached "datetimeFigureOut" date placeholder text (5/1/2025 -> 5/7/2025)
#    on every slide layout

$p = $ppt.ActivePresentation

# --- 1) Reposition the six repositioned pictures on slide 1 -----------------
$s = $p.Slides.Item(1)

# Shape.Id -> new Top (points). Values were tuned so that, after the engine's
# internal Single(float32)-precision round trip, the resulting EMU offset in
# the OOXML exactly matches the target (x stays untouched, only y/top moves).
$newTops = @{
    13 = 599.0200787401575   # 555426,7531612  -> 555426,7607555
    18 = 916.6112598425196   # 11921969,11565020 -> 11921969,11640963
    21 = 918.7770178740158   # 757340,11592525 -> 757340,11668468
    25 = 920.1603149606299   # 6084573,11610093 -> 6084573,11686036
    28 = 606.8954330708661   # 12046009,7631629 -> 12046009,7707572
    31 = 599.0200787401575   # 6084573,7531612 -> 6084573,7607555
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($newTops.ContainsKey($shape.Id)) {
        $shape.Top = $newTops[$shape.Id]
    }
}

# --- 2) Refresh the date placeholder text on every slide layout -------------
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "5/7/2025"
        }
    }
}
